# Add season-record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill in the season record for every data row (2-46)
$wins = 88
$losses = 74
$ties = 0

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
